$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new values in L1:R1 ---
$ws.Range("L1").Value = 0
$ws.Range("M1").Value = 0.00053272499999999995
$ws.Range("N1").Value = 0.001065
$ws.Range("O1").Value = 0.001598
$ws.Range("P1").Value = 0.0021310000000000001
$ws.Range("Q1").Value = 0.0042620000000000002
$ws.Range("R1").Value = 0.0063930000000000002

# --- Row 2: new difference formulas ---
$ws.Range("M2").Formula = "=N1-M1"
$ws.Range("N2").Formula = "=O1-N1"
$ws.Range("O2").Formula = "=P1-O1"
$ws.Range("P2").Formula = "=Q1-P1"
$ws.Range("Q2").Formula = "=R1-Q1"

# --- Row 5: k_min row (sets shared-string slot 31, reused from "Price:") ---
$ws.Range("G5").Value = "k_min"

# --- Row 4: header relabeling for the debug/calc block ---
$ws.Range("H4").Value = "Calc"
$ws.Range("I4").Value = "MX"

# --- Row 6: k_max row ---
$ws.Range("G6").Value = "k_max"

# --- Row 4 (cont.): error label ---
$ws.Range("J4").Value = "error"
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""

# --- Row 5 (cont.): numeric values and formula ---
$ws.Range("H5").Value = -0.041106482744587503
$ws.Range("I5").Value = -0.041106446492200001
$ws.Range("J5").Formula = "=I5-H5"

# --- Row 6 (cont.): numeric values and formula ---
$ws.Range("H6").Value = 0.045297248358925103
$ws.Range("I6").Value = 0.045297252020500001
$ws.Range("J6").Formula = "=I6-H6"

# --- Row 7: length row ---
$ws.Range("G7").Value = "length"
$ws.Range("H7").Formula = "=H6-H5"
$ws.Range("I7").Formula = "=I6-I5"

# --- Row 8: dk row ---
$ws.Range("G8").Value = "dk"
$ws.Range("H8").Value = 0.00090951295898434399

# --- Row 9: Nk row ---
$ws.Range("G9").Value = "Nk"
$ws.Range("H9").Formula = "=H7/H8"

# --- View adjustments ---
$ws.Columns("J").ColumnWidth = 11.140625
$win = $wb.Windows.Item(1)
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("M2").Select()
